# "added OPQA-3196 to enwiam052" -- the row for TCID "ENW000011"
# (Jira id "OPQA-3196", Description "As a user, I want to be able to
# see all emails that are associated to my Neon identity under the
# account page") is removed from the Test Cases sheet. Deleting the
# whole row shifts every row below it up by one and Excel automatically
# drops the now-unused shared strings from the string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(30).Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("D29").Select()
